# Pathway enrichment analysis update: the underlying statistics were
# recomputed (common_utils change) such that every p-value / padjust / FDR
# figure in both result tables is exactly twice its previous value, with
# the usual statistical cap of 1.0 (p-values/adjusted p-values cannot
# exceed 1).
#
# Sheet "A. mRNA-protein KEGG pathways": columns B (pvalue before),
#   D (padjust before), E (FDR before) -- rows 2..34
# Sheet "B. Residuals KEGG pathways": columns B (pvalue after),
#   C (padjust after), D (FDR after) -- rows 2..34

$wb = $excel.ActiveWorkbook

function Update-Column($ws, [int]$col, [int]$firstRow, [int]$lastRow) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $old = $cell.Value2
        if ($old -eq $null) { continue }
        $new = $old * 2
        if ($new -gt 1) { $new = 1 }
        $cell.Value2 = $new
    }
}

$ws1 = $wb.Worksheets.Item(1)
Update-Column $ws1 2 2 34   # column B
Update-Column $ws1 4 2 34   # column D
Update-Column $ws1 5 2 34   # column E

$ws2 = $wb.Worksheets.Item(2)
Update-Column $ws2 2 2 34   # column B
Update-Column $ws2 3 2 34   # column C
Update-Column $ws2 4 2 34   # column D
